# Fix group import: replace the SINH_VIEN (student code) values for groups
# 1-8 (rows 2-41) with the corrected "23211DH00xx" codes, restoring the
# previously-removed extra codes back to the top of the list. Group 9
# (rows 42-46) keeps its original student codes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCodes = @(
    "23211DH0001", "23211DH0002", "23211DH0003", "23211DH0004", "23211DH0005",
    "23211DH0006", "23211DH0007", "23211DH0008", "23211DH0009", "23211DH0010",
    "23211DH0011", "23211DH0012", "23211DH0013", "23211DH0014", "23211DH0015",
    "23211DH0016", "23211DH0017", "23211DH0018", "23211DH0019", "23211DH0020",
    "23211DH0021", "23211DH0022", "23211DH0023", "23211DH0024", "23211DH0025",
    "23211DH0026", "23211DH0027", "23211DH0028", "23211DH0029", "23211DH0030",
    "23211DH0031", "23211DH0032", "23211DH0033", "23211DH0034", "23211DH0035",
    "23211DH0036", "23211DH0037", "23211DH0038", "23211DH0039", "23211DH0040"
)

for ($i = 0; $i -lt $newCodes.Length; $i++) {
    $row = $i + 2
    $ws.Range("D$row").Value = $newCodes[$i]
}

# Restore the active selection to H7 (was J7 before the fix).
$ws.Range("H7").Select()
